$d = $word.ActiveDocument

# --- Step 1: turn "Bill Watson" into "William Watson" then split it into the
#     four runs ("W" | "ill" | "iam" | " Watson") that a live typing edit
#     (replace "B" with "W", leave "ill", insert "iam", leave " Watson")
#     would have produced.
$nameRange = $d.Range(0, 4)
$nameRange.Text = "William"

# Split the (currently single) run into separate runs at the boundaries
# between "W" / "ill" / "iam" / " Watson" by dropping a transient bookmark at
# each boundary and immediately deleting it again - Word leaves the run
# break behind even after the bookmark that caused it is gone.
foreach ($pos in 1, 4, 7) {
    $splitPoint = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TempSplitMarker", $splitPoint)
    $d.Bookmarks("TempSplitMarker").Delete()
}

# --- Step 2: move the "_GoBack" bookmark to the end of the (now longer)
#     name paragraph - this is where Word leaves it after the last edit.
$namePara = $d.Paragraphs(1)
$lastEditPos = $namePara.Range.End - 1
$goBackRange = $d.Range($lastEditPos, $lastEditPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
